$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# The textbox currently holds the text split across multiple runs
# ("The" / " " / "picture" / " " / "first"). Re-assigning the exact
# same text is treated as a no-op by the run-diffing text setter, so
# first set an unrelated placeholder to clear the old runs, then set
# the final desired text so it collapses into a single run.
$shp.TextFrame.TextRange.Text = "placeholder"
$shp.TextFrame.TextRange.Text = "The picture first"
